$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update day-hour values (Do column = E) for rows 2-7
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 5.5

$ws.Range("E3").Value = 6

$ws.Range("E4").Value = 6

$ws.Range("E5").Value = 6

$ws.Range("E6").Value = 6

$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = 0.5

# Update the active selection to match the target view
$ws.Range("E18").Select()
